$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters A..T correspond to the 20 fields of the LR-pairs table
$rowsData = @(
    ,("ECs", "Ucn2", "Il10rb", "ECs", 2, 0.6666666666666666, 0.1458113333333333, 0.437434, 0.1510638746080812, 0.1510638746080812, 3, 1, 81.038515, 243.115545, 0.8313546010859307, 0.8313546010859306, 11.81633392350333, 106.34700531153, 0.1255876472132964, 0.1255876472132964)
    ,("ECs", "Ucn2", "Il10rb", "FAPs", 2, 0.6666666666666666, 0.1458113333333333, 0.437434, 0.1510638746080812, 0.1510638746080812, 3, 1, 12.88200833333333, 38.646025, 0.1321534198787326, 0.1321534198787326, 1.878342811094445, 16.90508529985, 0.01996360764958997, 0.01996360764958997)
    ,("ECs", "Ucn2", "Il10rb", "sCs", 2, 0.6666666666666666, 0.1458113333333333, 0.437434, 0.1510638746080812, 0.1510638746080812, 3, 1, 3.557153333333333, 10.67146, 0.03649197903533674, 0.03649197903533673, 0.5186732704044446, 4.668059433640001, 0.005512619745194838, 0.005512619745194837)
    ,("FAPs", "Ucn2", "Il10rb", "ECs", 3, 1, 0.5810133333333333, 1.74304, 0.6019430954083812, 0.6019430954083812, 3, 1, 81.038515, 243.115545, 0.8313546010859307, 0.8313546010859306, 47.08445772853333, 423.7601195568, 0.5004281619596651, 0.500428161959665)
    ,("FAPs", "Ucn2", "Il10rb", "FAPs", 3, 1, 0.5810133333333333, 1.74304, 0.6019430954083812, 0.6019430954083812, 3, 1, 12.88200833333333, 38.646025, 0.1321534198787326, 0.1321534198787326, 7.484618601777777, 67.361567416, 0.07954883863060778, 0.07954883863060778)
    ,("FAPs", "Ucn2", "Il10rb", "sCs", 3, 1, 0.5810133333333333, 1.74304, 0.6019430954083812, 0.6019430954083812, 3, 1, 3.557153333333333, 10.67146, 0.03649197903533674, 0.03649197903533673, 2.066753515377778, 18.6007816384, 0.02196609481810835, 0.02196609481810835)
    ,("sCs", "Ucn2", "Il10rb", "ECs", 3, 1, 0.238405, 0.7152149999999999, 0.2469930299835376, 0.2469930299835376, 3, 1, 81.038515, 243.115545, 0.8313546010859307, 0.8313546010859306, 19.319987168575, 173.879884517175, 0.2053387919129692, 0.2053387919129692)
    ,("sCs", "Ucn2", "Il10rb", "FAPs", 3, 1, 0.238405, 0.7152149999999999, 0.2469930299835376, 0.2469930299835376, 3, 1, 12.88200833333333, 38.646025, 0.1321534198787326, 0.1321534198787326, 3.071135196708333, 27.640216770375, 0.03264097359853483, 0.03264097359853483)
    ,("sCs", "Ucn2", "Il10rb", "sCs", 3, 1, 0.238405, 0.7152149999999999, 0.2469930299835376, 0.2469930299835376, 3, 1, 3.557153333333333, 10.67146, 0.03649197903533674, 0.03649197903533673, 0.8480431404333333, 7.632388263899999, 0.009013264472033553, 0.009013264472033551)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
